# Auto-generated edit script applying the Hades_Profits.xlsx diff
# Updates currentAveragePrice/currentAveragePriceNQ/currentAveragePriceHQ/
# LevePriceNQ/LevePriceHQ/LeveProfitNQ/LeveProfitHQ values (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 74527.5
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 89333
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 267999
$ws.Range("M6").Value = -1388
$ws.Range("N6").Value = -268223
$ws.Range("H137").Value = 3227463.5
$ws.Range("I137").Value = 5883422
$ws.Range("K137").Value = 17650266
$ws.Range("M137").Value = -17647716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19013.96
$ws.Range("I32").Value = 17332.922
$ws.Range("J32").Value = 27979.5
$ws.Range("K32").Value = 17332.922
$ws.Range("L32").Value = 27979.5
$ws.Range("M32").Value = -17045.922
$ws.Range("N32").Value = -28553.5
$ws.Range("H74").Value = 12000643
$ws.Range("I74").Value = 13211163
$ws.Range("J74").Value = 500700
$ws.Range("K74").Value = 13211163
$ws.Range("L74").Value = 500700
$ws.Range("M74").Value = -13210289
$ws.Range("N74").Value = -502448
$ws.Range("H77").Value = 12000643
$ws.Range("I77").Value = 13211163
$ws.Range("J77").Value = 500700
$ws.Range("K77").Value = 66055815
$ws.Range("L77").Value = 2503500
$ws.Range("M77").Value = -66051447
$ws.Range("N77").Value = -2512236
$ws.Range("H122").Value = 1960.2222
$ws.Range("I122").Value = 1960.2222
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5880.6666
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3430.6666
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2478.261
$ws.Range("I134").Value = 2542.8572
$ws.Range("K134").Value = 7628.571599999999
$ws.Range("M134").Value = -5093.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1185.6666
$ws.Range("I16").Value = 719.3333
$ws.Range("J16").Value = 1885.1666
$ws.Range("K16").Value = 719.3333
$ws.Range("L16").Value = 1885.1666
$ws.Range("M16").Value = -432.3333
$ws.Range("N16").Value = -2459.1666
$ws.Range("H20").Value = 51000
$ws.Range("J20").Value = 51000
$ws.Range("L20").Value = 51000
$ws.Range("N20").Value = -51472
$ws.Range("H30").Value = 51000
$ws.Range("J30").Value = 51000
$ws.Range("L30").Value = 51000
$ws.Range("N30").Value = -51182
$ws.Range("H81").Value = 40328
$ws.Range("J81").Value = 40328
$ws.Range("L81").Value = 40328
$ws.Range("N81").Value = -42324
$ws.Range("H84").Value = 40328
$ws.Range("J84").Value = 40328
$ws.Range("L84").Value = 120984
$ws.Range("N84").Value = -130968
$ws.Range("H107").Value = 289.07407
$ws.Range("I107").Value = 218.65
$ws.Range("J107").Value = 490.2857
$ws.Range("K107").Value = 218.65
$ws.Range("L107").Value = 490.2857
$ws.Range("M107").Value = 1701.35
$ws.Range("N107").Value = -4330.2857
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 1185.6666
$ws.Range("I113").Value = 719.3333
$ws.Range("J113").Value = 1885.1666
$ws.Range("K113").Value = 719.3333
$ws.Range("L113").Value = 1885.1666
$ws.Range("M113").Value = 1450.6667
$ws.Range("N113").Value = -6225.1666
$ws.Range("H122").Value = 1400.5834
$ws.Range("I122").Value = 1127.909
$ws.Range("J122").Value = 4400
$ws.Range("K122").Value = 3383.727
$ws.Range("L122").Value = 13200
$ws.Range("M122").Value = -933.7270000000003
$ws.Range("N122").Value = -18100
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H128").Value = 51000
$ws.Range("J128").Value = 51000
$ws.Range("L128").Value = 51000
$ws.Range("N128").Value = -60960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1680269
$ws.Range("I11").Value = 3604291.8
$ws.Range("J11").Value = 455890.9
$ws.Range("K11").Value = 10812875.4
$ws.Range("L11").Value = 1367672.7
$ws.Range("M11").Value = -10812735.4
$ws.Range("N11").Value = -1367952.7
$ws.Range("H23").Value = 143.90909
$ws.Range("J23").Value = 138.11111
$ws.Range("L23").Value = 414.33333
$ws.Range("N23").Value = -884.3333299999999
$ws.Range("H64").Value = 3441.4482
$ws.Range("I64").Value = 1333.3334
$ws.Range("J64").Value = 3684.6924
$ws.Range("K64").Value = 4000.0002
$ws.Range("L64").Value = 11054.0772
$ws.Range("M64").Value = -3730.0002
$ws.Range("N64").Value = -11594.0772
$ws.Range("H67").Value = 3441.4482
$ws.Range("I67").Value = 1333.3334
$ws.Range("J67").Value = 3684.6924
$ws.Range("K67").Value = 4000.0002
$ws.Range("L67").Value = 11054.0772
$ws.Range("M67").Value = -3064.0002
$ws.Range("N67").Value = -12926.0772
$ws.Range("H87").Value = 9994
$ws.Range("I87").Value = 9994
$ws.Range("K87").Value = 29982
$ws.Range("M87").Value = -28734
$ws.Range("H90").Value = 9994
$ws.Range("I90").Value = 9994
$ws.Range("K90").Value = 89946
$ws.Range("M90").Value = -83706
$ws.Range("H114").Value = 15152175
$ws.Range("J114").Value = 20833900
$ws.Range("L114").Value = 62501700
$ws.Range("N114").Value = -62508208
$ws.Range("H119").Value = 4418.222
$ws.Range("I119").Value = 1350
$ws.Range("J119").Value = 6872.8
$ws.Range("K119").Value = 4050
$ws.Range("L119").Value = 20618.4
$ws.Range("M119").Value = 788
$ws.Range("N119").Value = -30294.4
$ws.Range("H120").Value = 11555.25
$ws.Range("I120").Value = 6000
$ws.Range("J120").Value = 13407
$ws.Range("K120").Value = 18000
$ws.Range("L120").Value = 40221
$ws.Range("M120").Value = -13162
$ws.Range("N120").Value = -49897
$ws.Range("H122").Value = 649.8095
$ws.Range("I122").Value = 310.92856
$ws.Range("J122").Value = 1327.5714
$ws.Range("K122").Value = 2798.35704
$ws.Range("L122").Value = 11948.1426
$ws.Range("M122").Value = -348.3570399999999
$ws.Range("N122").Value = -16848.1426
$ws.Range("H125").Value = 4481.25
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 4835.7144
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 14507.1432
$ws.Range("M125").Value = -1080
$ws.Range("N125").Value = -24347.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 11999.167
$ws.Range("J5").Value = 11999.167
$ws.Range("L5").Value = 11999.167
$ws.Range("N5").Value = -12223.167
$ws.Range("H132").Value = 85305.086
$ws.Range("I132").Value = 100942.6
$ws.Range("J132").Value = 74135.42999999999
$ws.Range("K132").Value = 302827.8
$ws.Range("L132").Value = 222406.29
$ws.Range("M132").Value = -300297.8
$ws.Range("N132").Value = -227466.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2640.111
$ws.Range("I122").Value = 2578.6155
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 7735.8465
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -5285.8465
$ws.Range("N122").Value = -13300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H136").Value = 54639.605
$ws.Range("I136").Value = 31524.848
$ws.Range("J136").Value = 207197
$ws.Range("K136").Value = 94574.54400000001
$ws.Range("L136").Value = 621591
$ws.Range("M136").Value = -92024.54400000001
$ws.Range("N136").Value = -626691
